$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Appendix1")
$ws2 = $wb.Worksheets.Item("Appendix2")

# --- Appendix1 sheet ---
# Row 2 person changed: Kieu Quoc Tuan / HE130002 / FPTUHCM -> Tran Thi Thuy Nguyen / He130020 / FPTUHN
$ws1.Range("B2").Value = "Trần Thị Thúy Nguyên"
$ws1.Range("C2").Value = "He130020"
$ws1.Range("D2").Value = "FPTUHN"
$ws1.Range("E2").Value = 135
$ws1.Range("F2").Value = 123

# Row 3 person text changed (shared string content updated): Nguyen Van A / fe001 -> Tran Thi Thuy Nguyen / He130020
$ws1.Range("B3").Value = "Trần Thị Thúy Nguyên"
$ws1.Range("C3").Value = "He130020"
$ws1.Range("E3").Value = 123
$ws1.Range("F3").Value = 321

# --- Appendix2 sheet ---
# Row 2 person text changed (shared string content updated): Nguyen Van A / fe001 -> Tran Thi Thuy Nguyen / He130020
$ws2.Range("B2").Value = "Trần Thị Thúy Nguyên"
$ws2.Range("C2").Value = "He130020"
$ws2.Range("E2").Value = 2000000

# Row 3 person changed: Kieu Quoc Tuan / HE130002 / FPTUHCM -> Tran Thi Thuy Nguyen / He130020 / FPTUHN
$ws2.Range("B3").Value = "Trần Thị Thúy Nguyên"
$ws2.Range("C3").Value = "He130020"
$ws2.Range("D3").Value = "FPTUHN"
$ws2.Range("E3").Value = 1000000
